$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 36

# Copy the header style from H1 into the new I1:J1 headers, then set their text
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill data rows 2-36: I = 1 (constant), J = same value as H (IP)
for ($r = 2; $r -le $lastRow; $r++) {
    $ipValue = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ipValue
}
